$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ==================================================================
# 1) FORMATTING FIRST
#    Copy formats from already-styled cells onto the new cells so
#    styles.xml keeps exactly its original 5 cellXfs (the diff shows
#    no change at all to styles.xml).
# ==================================================================

# B1 ("pacote") takes the col-E style (s=1); I1 ("comando...") takes the col-I style (s=2)
$ws.Range("E2").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("I2").Copy()
$ws.Range("I1").PasteSpecial(-4122)

# Rows 12-21 (bodytrack, fluidanimate) already exist; only B & C are new, both wrap-style (s=4)
$ws.Range("B2:C2").Copy()
$ws.Range("B12:C21").PasteSpecial(-4122)

# Rows 22-36 (freqmine, splash2, splash2x) are brand-new rows: A/C/D/E/I copy straight from row 2,
# then B is overwritten with the col-E style (s=1) used for the package name in this block.
$ws.Range("A2:E2").Copy()
$ws.Range("A22:E36").PasteSpecial(-4122)
$ws.Range("I2").Copy()
$ws.Range("I22:I36").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("B22:B26").PasteSpecial(-4122)
$ws.Range("B27:B31").PasteSpecial(-4122)
$ws.Range("B32:B36").PasteSpecial(-4122)

# Row 39 note cell: col-I style (s=2)
$ws.Range("I2").Copy()
$ws.Range("I39").PasteSpecial(-4122)

# ==================================================================
# 2) VALUES & FORMULAS
# ==================================================================

$ws.Range("B1").Value = "pacote"
$ws.Range("I1").Value = "comando de execução do benchmark"

# blackscholes

# vips

# bodytrack
$ws.Range("B12").Value = "bodytrack"
$ws.Range("C12").Value = "test"
$ws.Range("B13").Value = "bodytrack"
$ws.Range("C13").Value = "simdev"
$ws.Range("B14").Value = "bodytrack"
$ws.Range("C14").Value = "simsmall"
$ws.Range("B15").Value = "bodytrack"
$ws.Range("C15").Value = "simlarge"
$ws.Range("B16").Value = "bodytrack"
$ws.Range("C16").Value = "native"

# fluidanimate
$ws.Range("B17").Value = "fluidanimate"
$ws.Range("C17").Value = "test"
$ws.Range("B18").Value = "fluidanimate"
$ws.Range("C18").Value = "simdev"
$ws.Range("B19").Value = "fluidanimate"
$ws.Range("C19").Value = "simsmall"
$ws.Range("B20").Value = "fluidanimate"
$ws.Range("C20").Value = "simlarge"
$ws.Range("B21").Value = "fluidanimate"
$ws.Range("C21").Value = "native"

# freqmine
$ws.Range("A22").Value = "021"
$ws.Range("B22").Value = "freqmine"
$ws.Range("C22").Value = "test"
$ws.Range("D22").Value = "parsecmgmt -a run"
$ws.Range("E22").Value = "result/exec-"
$ws.Range("I22").Formula = '=CONCATENATE(D22," -p ", B22," -i ",C22, "  >  ",E22,F22,A22,"-",B22,"-",C22,".txt")'
$ws.Range("A23").Value = "022"
$ws.Range("B23").Value = "freqmine"
$ws.Range("C23").Value = "simdev"
$ws.Range("D23").Value = "parsecmgmt -a run"
$ws.Range("E23").Value = "result/exec-"
$ws.Range("I23").Formula = '=CONCATENATE(D23," -p ", B23," -i ",C23, "  >  ",E23,F23,A23,"-",B23,"-",C23,".txt")'
$ws.Range("A24").Value = "023"
$ws.Range("B24").Value = "freqmine"
$ws.Range("C24").Value = "simsmall"
$ws.Range("D24").Value = "parsecmgmt -a run"
$ws.Range("E24").Value = "result/exec-"
$ws.Range("I24").Formula = '=CONCATENATE(D24," -p ", B24," -i ",C24, "  >  ",E24,F24,A24,"-",B24,"-",C24,".txt")'
$ws.Range("A25").Value = "024"
$ws.Range("B25").Value = "freqmine"
$ws.Range("C25").Value = "simlarge"
$ws.Range("D25").Value = "parsecmgmt -a run"
$ws.Range("E25").Value = "result/exec-"
$ws.Range("I25").Formula = '=CONCATENATE(D25," -p ", B25," -i ",C25, "  >  ",E25,F25,A25,"-",B25,"-",C25,".txt")'
$ws.Range("A26").Value = "025"
$ws.Range("B26").Value = "freqmine"
$ws.Range("C26").Value = "native"
$ws.Range("D26").Value = "parsecmgmt -a run"
$ws.Range("E26").Value = "result/exec-"
$ws.Range("I26").Formula = '=CONCATENATE(D26," -p ", B26," -i ",C26, "  >  ",E26,F26,A26,"-",B26,"-",C26,".txt")'

# splash2
$ws.Range("A27").Value = "026"
$ws.Range("B27").Value = "splash2"
$ws.Range("C27").Value = "test"
$ws.Range("D27").Value = "parsecmgmt -a run"
$ws.Range("E27").Value = "result/exec-"
$ws.Range("I27").Formula = '=CONCATENATE(D27," -p ", B27," -i ",C27, "  >  ",E27,F27,A27,"-",B27,"-",C27,".txt")'
$ws.Range("A28").Value = "027"
$ws.Range("B28").Value = "splash2"
$ws.Range("C28").Value = "simdev"
$ws.Range("D28").Value = "parsecmgmt -a run"
$ws.Range("E28").Value = "result/exec-"
$ws.Range("I28").Formula = '=CONCATENATE(D28," -p ", B28," -i ",C28, "  >  ",E28,F28,A28,"-",B28,"-",C28,".txt")'
$ws.Range("A29").Value = "028"
$ws.Range("B29").Value = "splash2"
$ws.Range("C29").Value = "simsmall"
$ws.Range("D29").Value = "parsecmgmt -a run"
$ws.Range("E29").Value = "result/exec-"
$ws.Range("I29").Formula = '=CONCATENATE(D29," -p ", B29," -i ",C29, "  >  ",E29,F29,A29,"-",B29,"-",C29,".txt")'
$ws.Range("A30").Value = "029"
$ws.Range("B30").Value = "splash2"
$ws.Range("C30").Value = "simlarge"
$ws.Range("D30").Value = "parsecmgmt -a run"
$ws.Range("E30").Value = "result/exec-"
$ws.Range("I30").Formula = '=CONCATENATE(D30," -p ", B30," -i ",C30, "  >  ",E30,F30,A30,"-",B30,"-",C30,".txt")'
$ws.Range("A31").Value = "030"
$ws.Range("B31").Value = "splash2"
$ws.Range("C31").Value = "native"
$ws.Range("D31").Value = "parsecmgmt -a run"
$ws.Range("E31").Value = "result/exec-"
$ws.Range("I31").Formula = '=CONCATENATE(D31," -p ", B31," -i ",C31, "  >  ",E31,F31,A31,"-",B31,"-",C31,".txt")'

# splash2x
$ws.Range("A32").Value = "031"
$ws.Range("B32").Value = "splash2x"
$ws.Range("C32").Value = "test"
$ws.Range("D32").Value = "parsecmgmt -a run"
$ws.Range("E32").Value = "result/exec-"
$ws.Range("I32").Formula = '=CONCATENATE(D32," -p ", B32," -i ",C32, "  >  ",E32,F32,A32,"-",B32,"-",C32,".txt")'
$ws.Range("A33").Value = "032"
$ws.Range("B33").Value = "splash2x"
$ws.Range("C33").Value = "simdev"
$ws.Range("D33").Value = "parsecmgmt -a run"
$ws.Range("E33").Value = "result/exec-"
$ws.Range("I33").Formula = '=CONCATENATE(D33," -p ", B33," -i ",C33, "  >  ",E33,F33,A33,"-",B33,"-",C33,".txt")'
$ws.Range("A34").Value = "033"
$ws.Range("B34").Value = "splash2x"
$ws.Range("C34").Value = "simsmall"
$ws.Range("D34").Value = "parsecmgmt -a run"
$ws.Range("E34").Value = "result/exec-"
$ws.Range("I34").Formula = '=CONCATENATE(D34," -p ", B34," -i ",C34, "  >  ",E34,F34,A34,"-",B34,"-",C34,".txt")'
$ws.Range("A35").Value = "034"
$ws.Range("B35").Value = "splash2x"
$ws.Range("C35").Value = "simlarge"
$ws.Range("D35").Value = "parsecmgmt -a run"
$ws.Range("E35").Value = "result/exec-"
$ws.Range("I35").Formula = '=CONCATENATE(D35," -p ", B35," -i ",C35, "  >  ",E35,F35,A35,"-",B35,"-",C35,".txt")'
$ws.Range("A36").Value = "035"
$ws.Range("B36").Value = "splash2x"
$ws.Range("C36").Value = "native"
$ws.Range("D36").Value = "parsecmgmt -a run"
$ws.Range("E36").Value = "result/exec-"
$ws.Range("I36").Formula = '=CONCATENATE(D36," -p ", B36," -i ",C36, "  >  ",E36,F36,A36,"-",B36,"-",C36,".txt")'

# Row 39 footnote (plain text, not a formula)
$ws.Range("I39").Value = "parsecmgmt -a run -p splash2x.fft -i simsmall -n 4  >  result/exec-200-splash2x-fft-simsmall.txt"

# ==================================================================
# 3) I3:I21 drop their t="shared" formula grouping in the target --
#    re-enter each as an independent (textually identical) formula.
# ==================================================================
$ws.Range("I3").Formula = '=CONCATENATE(D3," -p ", B3," -i ",C3, "  >  ",E3,F3,A3,"-",B3,"-",C3,".txt")'
$ws.Range("I4").Formula = '=CONCATENATE(D4," -p ", B4," -i ",C4, "  >  ",E4,F4,A4,"-",B4,"-",C4,".txt")'
$ws.Range("I5").Formula = '=CONCATENATE(D5," -p ", B5," -i ",C5, "  >  ",E5,F5,A5,"-",B5,"-",C5,".txt")'
$ws.Range("I6").Formula = '=CONCATENATE(D6," -p ", B6," -i ",C6, "  >  ",E6,F6,A6,"-",B6,"-",C6,".txt")'
$ws.Range("I7").Formula = '=CONCATENATE(D7," -p ", B7," -i ",C7, "  >  ",E7,F7,A7,"-",B7,"-",C7,".txt")'
$ws.Range("I8").Formula = '=CONCATENATE(D8," -p ", B8," -i ",C8, "  >  ",E8,F8,A8,"-",B8,"-",C8,".txt")'
$ws.Range("I9").Formula = '=CONCATENATE(D9," -p ", B9," -i ",C9, "  >  ",E9,F9,A9,"-",B9,"-",C9,".txt")'
$ws.Range("I10").Formula = '=CONCATENATE(D10," -p ", B10," -i ",C10, "  >  ",E10,F10,A10,"-",B10,"-",C10,".txt")'
$ws.Range("I11").Formula = '=CONCATENATE(D11," -p ", B11," -i ",C11, "  >  ",E11,F11,A11,"-",B11,"-",C11,".txt")'
$ws.Range("I12").Formula = '=CONCATENATE(D12," -p ", B12," -i ",C12, "  >  ",E12,F12,A12,"-",B12,"-",C12,".txt")'
$ws.Range("I13").Formula = '=CONCATENATE(D13," -p ", B13," -i ",C13, "  >  ",E13,F13,A13,"-",B13,"-",C13,".txt")'
$ws.Range("I14").Formula = '=CONCATENATE(D14," -p ", B14," -i ",C14, "  >  ",E14,F14,A14,"-",B14,"-",C14,".txt")'
$ws.Range("I15").Formula = '=CONCATENATE(D15," -p ", B15," -i ",C15, "  >  ",E15,F15,A15,"-",B15,"-",C15,".txt")'
$ws.Range("I16").Formula = '=CONCATENATE(D16," -p ", B16," -i ",C16, "  >  ",E16,F16,A16,"-",B16,"-",C16,".txt")'
$ws.Range("I17").Formula = '=CONCATENATE(D17," -p ", B17," -i ",C17, "  >  ",E17,F17,A17,"-",B17,"-",C17,".txt")'
$ws.Range("I18").Formula = '=CONCATENATE(D18," -p ", B18," -i ",C18, "  >  ",E18,F18,A18,"-",B18,"-",C18,".txt")'
$ws.Range("I19").Formula = '=CONCATENATE(D19," -p ", B19," -i ",C19, "  >  ",E19,F19,A19,"-",B19,"-",C19,".txt")'
$ws.Range("I20").Formula = '=CONCATENATE(D20," -p ", B20," -i ",C20, "  >  ",E20,F20,A20,"-",B20,"-",C20,".txt")'
$ws.Range("I21").Formula = '=CONCATENATE(D21," -p ", B21," -i ",C21, "  >  ",E21,F21,A21,"-",B21,"-",C21,".txt")'

# ==================================================================
# 4) Column width + sheet view (scrolled so row 22 is at top, with the
#    splash2 block A27:I36 selected -- matches the saved view state).
# ==================================================================
$ws.Columns("I").ColumnWidth = 79.85546875
$ws.Application.Goto($ws.Range("A22"))
$ws.Range("A27:I36").Select()
